$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy only the formatting (borders, number formats, alignment) from the last
# existing data row (19) down into the new row (20).
$ws.Range("B19:G19").Copy()
$ws.Range("B20:G20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New weekly backup entry (one week after row 19 -> 2015-07-10)
$ws.Range("B20").Value = 42195
$ws.Range("C20").Value = "Jovanny Zepeda"
$ws.Range("E20").Value = 42195
$ws.Range("G20").Value = "Realizada"

# Re-merge the split columns for the new row
[void]$ws.Range("C20:D20").Merge()
[void]$ws.Range("E20:F20").Merge()

# Move the active selection down to the next empty row, as in the saved file
[void]$ws.Range("E21").Select()
